# Automatische test-sync: 2025-06-26 18:48:50
# Adds a "Handmatig opvolgen" column + a new logged e-mail row to the Logs
# sheet, mirrors the new category count on the Dashboard sheet, and widens
# the chart series + conditional-formatting ranges to include the new row.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: new "Handmatig opvolgen" header column ---------------
$logs.Range("G1").Copy() | Out-Null
$logs.Range("H1").PasteSpecial(-4122) | Out-Null
$logs.Range("H1").Value = "Handmatig opvolgen"

# --- Logs sheet: new data row (row 3) ----------------------------------
$logs.Range("A3").Value = "Kunnen jullie ook maatwerk leveren voor beurzen?"
$logs.Range("B3").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Range("C3").Value = "Goedemiddag, `nWij zijn op zoek naar maatwerk displays voor een beurs in september. Kunnen jullie dit soort werk leveren, en zo ja, wat zijn de mogelijkheden en levertijden?`nAlvast bedankt!`nMet vriendelijke groet, `nM. Peters`nMarketing & Events`nSent using {0}"
$logs.Range("D3").Value = "Offerte / Prijsaanvraag"
$logs.Range("E3").Value = "Beste M. Peters,`nHartelijk dank voor uw interesse in maatwerk displays voor uw beurs in september. Ja, wij kunnen zeker maatwerk displays leveren. Om u goed te kunnen helpen, ontvangen we graag meer informatie over uw specifieke wensen en eisen voor de displays. Hierdoor kunnen we u voorzien van de best passende opties en bijbehorende levertijden.`nKunt u meer details geven over de gewenste afmetingen, materialen, eventuele speciale functies, en de hoeveelheid displays die u nodig heeft? Met deze informatie kunnen we een passende offerte op maat voor u opstellen en uitleggen wat de mogelijkheden zijn.`nIndien u ons meer details kunt verstrekken, kunnen we samen verder kijken naar de beste oplossing voor uw beurs in september.`nMet vriendelijke groet,`n[Bedrijfsnaam] `nE-mailassistent"
$logs.Range("F3").Value = "2025-06-26 18:47:57"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Nee"

# Row 3 holds multi-line cell text (C3/E3); clear the implicit autofit height
# the engine stamps on write so the row stays unsized, like row 2.
$logs.Rows.Item(3).AutoFit()

# --- Logs sheet: extend existing conditional formatting down to row 3 --
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D3")) | Out-Null
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G3")) | Out-Null

# --- Logs sheet: new conditional formatting rule for column H ----------
$fc = $logs.Range("H2:H3").FormatConditions.Add(1, 3, '"Ja"')
$fc.Interior.Color = 13431551
$fc.Priority = 9

# --- Dashboard sheet: new category count row (row 3) -------------------
$dash.Range("A3").Value = "Offerte / Prijsaanvraag"
$dash.Range("B3").Value = 1

# --- Dashboard chart: widen category/value series to include row 3 -----
$chart = $dash.ChartObjects(1).Chart
$chart.SeriesCollection(1).XValues = "='Dashboard'!`$A`$2:`$A`$3"
$chart.SeriesCollection(1).Values = "='Dashboard'!`$B`$2:`$B`$3"
